$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 85.77505782882612

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
